$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.273.00"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.706.96"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.45"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.29"
$ws.Range("E6").Value = "  +4.46%  "
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.708.73"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.362"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.34"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.203.87"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.276.66"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.716.29"
$ws.Range("E18").Value = "  +2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.75"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "369.21"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "576.80"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("E34").Value = "  +5.55%  "
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.82"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "159.35"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.96"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0309"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.72"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.594"
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "154.79"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("E51").Value = "  +3.39%  "
